{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\ntable.getCell(0, 0).value = \"753\u00d78=6024\";\ntable.getCell(0, 1).value = \"202\u00d74=808\";\ntable.getCell(0, 2).value = \"845\u00d78=6760\";\ntable.getCell(0, 3).value = \"339\u00d78=2712\";\ntable.getCell(0, 4).value = \"402\u00d72=804\";\ntable.getCell(4, 0).value = \"445\u00d78=3560\";\ntable.getCell(4, 1).value = \"881\u00d74=3524\";\ntable.getCell(4, 2).value = \"136\u00d75=680\";\ntable.getCell(4, 3).value = \"676\u00d79=6084\";\ntable.getCell(4, 4).value = \"681\u00d76=4086\";\ntable.getCell(9, 0).value = \"336\u00d74=1344\";\ntable.getCell(9, 1).value = \"529\u00d74=2116\";\ntable.getCell(9, 2).value = \"903\u00d77=6321\";\ntable.getCell(9, 3).value = \"602\u00d76=3612\";\ntable.getCell(9, 4).value = \"768\u00d75=3840\";\ntable.getCell(14, 0).value = \"796\u00d77=5572\";\ntable.getCell(14, 1).value = \"878\u00d75=4390\";\ntable.getCell(14, 2).value = \"871\u00d76=5226\";\ntable.getCell(14, 3).value = \"542\u00d73=1626\";\ntable.getCell(14, 4).value = \"675\u00d78=5400\";\ntable.getCell(19, 0).value = \"466\u00d78=3728\";\ntable.getCell(19, 1).value = \"213\u00d73=639\";\ntable.getCell(19, 2).value = \"657\u00d73=1971\";\ntable.getCell(19, 3).value = \"896\u00d76=5376\";\ntable.getCell(19, 4).value = \"719\u00d79=6471\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$table.Cell(1,1).Range.Text = \"753\u00d78=6024\"\n$table.Cell(1,2).Range.Text = \"202\u00d74=808\"\n$table.Cell(1,3).Range.Text = \"845\u00d78=6760\"\n$table.Cell(1,4).Range.Text = \"339\u00d78=2712\"\n$table.Cell(1,5).Range.Text = \"402\u00d72=804\"\n$table.Cell(5,1).Range.Text = \"445\u00d78=3560\"\n$table.Cell(5,2).Range.Text = \"881\u00d74=3524\"\n$table.Cell(5,3).Range.Text = \"136\u00d75=680\"\n$table.Cell(5,4).Range.Text = \"676\u00d79=6084\"\n$table.Cell(5,5).Range.Text = \"681\u00d76=4086\"\n$table.Cell(10,1).Range.Text = \"336\u00d74=1344\"\n$table.Cell(10,2).Range.Text = \"529\u00d74=2116\"\n$table.Cell(10,3).Range.Text = \"903\u00d77=6321\"\n$table.Cell(10,4).Range.Text = \"602\u00d76=3612\"\n$table.Cell(10,5).Range.Text = \"768\u00d75=3840\"\n$table.Cell(15,1).Range.Text = \"796\u00d77=5572\"\n$table.Cell(15,2).Range.Text = \"878\u00d75=4390\"\n$table.Cell(15,3).Range.Text = \"871\u00d76=5226\"\n$table.Cell(15,4).Range.Text = \"542\u00d73=1626\"\n$table.Cell(15,5).Range.Text = \"675\u00d78=5400\"\n$table.Cell(20,1).Range.Text = \"466\u00d78=3728\"\n$table.Cell(20,2).Range.Text = \"213\u00d73=639\"\n$table.Cell(20,3).Range.Text = \"657\u00d73=1971\"\n$table.Cell(20,4).Range.Text = \"896\u00d76=5376\"\n$table.Cell(20,5).Range.Text = \"719\u00d79=6471\"\n"}
